# Weekly update to the fruit/vegetable price sheet: a new daily record is
# inserted at row 33 (pushing the existing rows 33-93 down to 34-94), adding
# one more observation to the "Espárragos" series for
# "Feria Lagunitas de Puerto Montt".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 33; Excel shifts rows 33:93 down to 34:94
# and carries formatting (e.g. the date number format on column D) along.
$ws.Rows.Item(33).Insert()

$ws.Cells.Item(33, 1).Value = 4
$ws.Cells.Item(33, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(33, 3).Value = "Los Lagos"
$ws.Cells.Item(33, 4).Value = 45260
$ws.Cells.Item(33, 5).Value = 10
$ws.Cells.Item(33, 6).Value = 300000000
$ws.Cells.Item(33, 7).Value = "Espárragos"
$ws.Cells.Item(33, 8).Value = "Sin especificar"
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 200
$ws.Cells.Item(33, 11).Value = 2000
$ws.Cells.Item(33, 12).Value = 2000
$ws.Cells.Item(33, 13).Value = 2000
$ws.Cells.Item(33, 14).Value = "$/kilo"
$ws.Cells.Item(33, 15).Value = "Provincia de Linares"
$ws.Cells.Item(33, 16).Value = 2000
$ws.Cells.Item(33, 17).Value = 1
$ws.Cells.Item(33, 18).Value = "Hortaliza"
